# Automatic update of files.
# Update the "Förändrad" date column (C) for rows 2-15 from 45179 to 45180
# (i.e. bump the serial date value by one day, 2023-09-10 -> 2023-09-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
